$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(newPrice, newVolumePct)  (either may be $null to leave unchanged)
$updates = @{
    2 = @("59.839.10", "  +2.68%  ")
    3 = @("2.418.85", "  +3.11%  ")
    4 = @("0.999", "  -0.07%  ")
    5 = @("554.74", "  +2.22%  ")
    6 = @("137.46", "  +1.81%  ")
    7 = @($null, "  -0.12%  ")
    8 = @("0.568", "  +1.27%  ")
    9 = @($null, "  +5.33%  ")
    10 = @($null, "  +2.21%  ")
    11 = @("0.363", "  +1.91%  ")
    13 = @("24.64", $null)
    14 = @("2.842.90", "  +2.88%  ")
    15 = @("59.666.02", "  +2.50%  ")
    16 = @("0.0000140", "  +4.37%  ")
    17 = @("2.432.62", "  +3.72%  ")
    18 = @("11.33", "  +5.63%  ")
    19 = @($null, "  +4.65%  ")
    20 = @("336.67", "  +0.99%  ")
    21 = @("6.98", "  +4.69%  ")
    22 = @($null, "  +0.01%  ")
    23 = @("64.58", "  +2.84%  ")
    24 = @($null, "  +0.99%  ")
    25 = @("8.56", "  +0.20%  ")
    26 = @("1.00", "  -0.05%  ")
    27 = @("1.40", "  -1.15%  ")
    28 = @($null, "  +6.48%  ")
    29 = @("1.81", "  +2.77%  ")
    30 = @("170.54", "  +0.09%  ")
    31 = @("6.27", "  +2.63%  ")
    32 = @("18.72", "  +1.64%  ")
    33 = @("1.03", "  -0.17%  ")
    34 = @($null, "  -0.02%  ")
    35 = @("1.32", "  +5.58%  ")
    36 = @("4.31", "  +1.05%  ")
    37 = @("0.999", "  -0.04%  ")
    38 = @("1.65", "  +0.07%  ")
    39 = @("40.18", "  +2.65%  ")
    40 = @($null, "  +11.22%  ")
    41 = @("306.54", "  +6.41%  ")
    42 = @($null, "  +3.14%  ")
    43 = @("142.52", "  +0.02%  ")
    44 = @("0.0964", "  +2.88%  ")
    45 = @("0.0527", "  +4.60%  ")
    46 = @("0.574", "  +1.61%  ")
    47 = @("19.12", "  -0.32%  ")
    48 = @("0.407", "  +6.62%  ")
    49 = @($null, "  +3.40%  ")
    50 = @("11.04", "  -0.26%  ")
    51 = @("1.62", "  +5.63%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $pct = $pair[1]
    if ($null -ne $price) {
        $cell = $ws.Range("D$row")
        # force text storage so numeric-looking strings (e.g. "0.999", "1.00")
        # are not smart-converted to numbers, matching the source inlineStr cells
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.ClearFormats()
    }
    if ($null -ne $pct) {
        $ws.Range("E$row").Value = $pct
    }
}
